$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.226.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.92%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.933.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.96%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.07%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  +0.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.90"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.70%  "

$ws.Range("E10").Value = "  +0.18%  "

$ws.Range("E11").Value = "  -1.47%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000228"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.27%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.26%  "

$ws.Range("E14").Value = "  -0.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.420.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.199.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.97%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.49%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.928.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "432.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.43%  "

$ws.Range("E20").Value = "  -1.46%  "

$ws.Range("E21").Value = "  +1.70%  "

$ws.Range("E22").Value = "  -0.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.92%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("E25").Value = "  +0.94%  "

$ws.Range("E26").Value = "  +1.90%  "

$ws.Range("E27").Value = "  -0.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.57%  "

$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.00%  "

$ws.Range("E33").Value = "  +1.52%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0861"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.95%  "

$ws.Range("E35").Value = "  +0.61%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.10%  "

$ws.Range("E38").Value = "  +2.82%  "

$ws.Range("E39").Value = "  -1.53%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.74%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.289"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.79%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.34%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "377.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.82%  "

$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.718.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.42%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "130.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.82%  "

$ws.Range("E47").Value = "  +0.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.29%  "

$ws.Range("E49").Value = "  -0.04%  "

$ws.Range("E50").Value = "  -2.54%  "

$ws.Range("E51").Value = "  +2.21%  "
